$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New criterion/score rows replacing the old "Full Evaluation" row
$criteria = @(
    "*   **Coherence",
    "*   **Accuracy of Financial Insights",
    "*   **Value to Investors",
    "*   **Clarity of Writing",
    "*   **Coherence",
    "*   **Accuracy of Financial Insights",
    "*   **Value to Investors",
    "*   **Clarity of Writing"
)

$scores = @("8", "9", "9", "9", "8", "9", "9", "9")

for ($i = 0; $i -lt $criteria.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $criteria[$i]
    $ws.Cells.Item($row, 2).Value = "'" + $scores[$i]
}
